# Apply the "Cleaned data for weather and survey vars" update to the
# Data Dictionary workbook:
#   - refresh two FEMA source links (rows 15 & 17)
#   - append 7 new variable-definition rows (20-26) describing the newly
#     added survey demographic fields, household composition fields, the
#     survey wave/date field, and the case_ID row identifier
#   - leave the active selection on B20, matching the authored edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New rows describing the additional survey variables ---
# (values entered in the same row/column order the workbook was authored in)

# Row 25: survey wave / date
$ws.Range("A25").Value = "Time"
$ws.Range("B25").Value = "wave"
$ws.Range("C25").Value = "date of survey"
$ws.Range("D25").Value = "date"
$ws.Range("G25").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H25").Value = "https://osf.io/jw79p/"

# Row 20: religion
$ws.Range("A20").Value = "Demographic"
$ws.Range("B20").Value = "religion"
$ws.Range("G20").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H20").Value = "https://osf.io/jw79p/"
$ws.Range("D20").Value = "cat"

# Row 21: marital status
$ws.Range("A21").Value = "Demographic"
$ws.Range("B21").Value = "marit_status"
$ws.Range("G21").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H21").Value = "https://osf.io/jw79p/"
$ws.Range("D21").Value = "cat"

# Row 22: employment
$ws.Range("A22").Value = "Demographic"
$ws.Range("B22").Value = "employment"
$ws.Range("G22").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H22").Value = "https://osf.io/jw79p/"
$ws.Range("D22").Value = "cat"

# Row 20-22 descriptions (column C)
$ws.Range("C20").Value = "stated religion"
$ws.Range("C21").Value = "stated marital status"
$ws.Range("C22").Value = "current employment"

# Row 23: children in household
$ws.Range("A23").Value = "Demographic"
$ws.Range("B23").Value = "children"
$ws.Range("G23").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H23").Value = "https://osf.io/jw79p/"
$ws.Range("D23").Value = "int"
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10

# Row 24: adults in household
$ws.Range("A24").Value = "Demographic"
$ws.Range("B24").Value = "adults"
$ws.Range("G24").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H24").Value = "https://osf.io/jw79p/"
$ws.Range("D24").Value = "int"
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 13

# Row 23-24 descriptions (column C)
$ws.Range("C23").Value = "number of children in household"
$ws.Range("C24").Value = "number of adults in household"

# Row 26: row/case identifier
$ws.Range("A26").Value = "Row Mane"
$ws.Range("B26").Value = "case_ID"
$ws.Range("C26").Value = "research ID to track row"
$ws.Range("D26").Value = "int"
$ws.Range("G26").Value = "Leiserowitz et. al. (2022)"
$ws.Range("H26").Value = "https://osf.io/jw79p/"

# --- Row 15 & 17: point the FEMA source links at the current OpenFEMA pages ---
$ws.Range("H17").Value = "https://www.fema.gov/openfema-data-page/public-assistance-funded-projects-details-v1"
$ws.Range("H15").Value = "https://www.fema.gov/openfema-data-page/disaster-declarations-summaries-v2"

# --- Match formatting used by the rest of the table ---
# Columns A, C & D of the existing rows use the bold "Arial 10" style (s=1);
# copy that formatting from row 5 (an existing fully-styled data row) onto
# the new rows without disturbing the shared string table / styles.xml.
# PasteSpecial only reliably honours the first contiguous block of a
# multi-area Union, so paste into each Areas member individually.
function Paste-FormatOnly($sourceRange, $targetRange) {
    $sourceRange.Copy()
    for ($i = 1; $i -le $targetRange.Areas.Count; $i++) {
        $targetRange.Areas.Item($i).PasteSpecial(-4122)
    }
}

Paste-FormatOnly $ws.Range("A5") $excel.Union($ws.Range("A20"), $ws.Range("A21"), $ws.Range("A22"), $ws.Range("A23"), $ws.Range("A24"), $ws.Range("A26"))
Paste-FormatOnly $ws.Range("C5") $excel.Union($ws.Range("C20"), $ws.Range("C21"), $ws.Range("C22"), $ws.Range("C23"), $ws.Range("C24"), $ws.Range("C25"), $ws.Range("C26"))
Paste-FormatOnly $ws.Range("D7") $excel.Union($ws.Range("D20"), $ws.Range("D21"), $ws.Range("D22"))
Paste-FormatOnly $ws.Range("D6") $excel.Union($ws.Range("D23"), $ws.Range("D24"), $ws.Range("D26"))

$excel.CutCopyMode = $false

# --- Final selection, matching the authored workbook state ---
[void]$ws.Range("B20").Select()
